$d = $word.ActiveDocument

# The commit removes two whole intro paragraphs:
#   1. "Another token game as Le Monde mathematical puzzle:"
#   2. "Which cannot be solved in a few lines of R code:"
# Walk the paragraphs back-to-front (so deleting one doesn't invalidate the
# index of paragraphs we still need to visit) and drop the Range (which
# includes the paragraph mark) for any paragraph matching either of those.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -like "*Which cannot be solved in a few lines of*R code*") {
        $p.Range.Delete()
    }
    elseif ($t -like "*Another token game as*Le Monde mathematical puzzle*") {
        $p.Range.Delete()
    }
}
